$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B (ID Competição) values were dropped / incorrectly scraped as 57
# for rows 2-34; recover the correct value of 257.
$ws.Range("B2:B34").Value = 257
